# Scheduled-runner update: refresh market-board price snapshots (and the
# dependent profit/loss columns) across the per-job "Phantom Profits" sheets.
# Values below are the latest pulled averages; a handful of rows also gain or
# lose a stray M/N profit cell where HQ/NQ pricing data became (un)available.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 350
$ws.Range("I2").Value = 350
$ws.Range("K2").Value = 350
$ws.Range("M2").Value = -237
$ws.Range("H17").Value = 2859.25
$ws.Range("J17").Value = 2859.25
$ws.Range("L17").Value = 8577.75
$ws.Range("N17").Value = -8913.75
$ws.Range("H64").Value = 4669.25
$ws.Range("I64").Value = 5265.75
$ws.Range("K64").Value = 5265.75
$ws.Range("M64").Value = -5017.75
$ws.Range("H67").Value = 4669.25
$ws.Range("I67").Value = 5265.75
$ws.Range("K67").Value = 5265.75
$ws.Range("M67").Value = -4407.75
$ws.Range("H127").Value = 3323.3333
$ws.Range("I127").Value = 2485
$ws.Range("K127").Value = 7455
$ws.Range("M127").Value = -2495
$ws.Range("H138").Value = 4850.227
$ws.Range("J138").Value = 4018.5334
$ws.Range("L138").Value = 12055.6002
$ws.Range("N138").Value = -22335.6002

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 6666
$ws.Range("I38").Value = 6666
$ws.Range("K38").Value = 6666
$ws.Range("M38").Value = -6199
$ws.Range("H43").Value = 38333.332
$ws.Range("I43").Value = 32500
$ws.Range("J43").Value = 50000
$ws.Range("K43").Value = 32500
$ws.Range("L43").Value = 50000
$ws.Range("M43").Value = -32187
$ws.Range("N43").Value = -50626
$ws.Range("H61").Value = 2782.8
$ws.Range("I61").Value = 2836.4443
$ws.Range("K61").Value = 2836.4443
$ws.Range("M61").Value = -2624.4443
$ws.Range("H103").Value = 40181
$ws.Range("J103").Value = 40181
$ws.Range("L103").Value = 40181
$ws.Range("N103").Value = -42525
$ws.Range("H122").Value = 3017.12
$ws.Range("I122").Value = 3187.6365
$ws.Range("K122").Value = 9562.9095
$ws.Range("M122").Value = -7112.9095
$ws.Range("H132").Value = 4430.727
$ws.Range("I132").Value = 5027.2354
$ws.Range("J132").Value = 2402.6
$ws.Range("K132").Value = 15081.7062
$ws.Range("L132").Value = 7207.799999999999
$ws.Range("M132").Value = -12551.7062
$ws.Range("N132").Value = -12267.8
$ws.Range("H136").Value = 2782.8
$ws.Range("I136").Value = 2836.4443
$ws.Range("K136").Value = 8509.332900000001
$ws.Range("M136").Value = -5959.332900000001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 23444
$ws.Range("I7").Value = 1999.5
$ws.Range("J7").Value = 34166.25
$ws.Range("K7").Value = 1999.5
$ws.Range("L7").Value = 34166.25
$ws.Range("M7").Value = -1886.5
$ws.Range("N7").Value = -34392.25
$ws.Range("H22").Value = 602.5
$ws.Range("I22").Value = 641.75
$ws.Range("K22").Value = 641.75
$ws.Range("M22").Value = -468.75
$ws.Range("H37").Value = 469
$ws.Range("I37").Value = 469
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 469
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -332
$ws.Range("N37").ClearContents()
$ws.Range("H38").Value = 35299.668
$ws.Range("I38").Value = 26000
$ws.Range("J38").Value = 39949.5
$ws.Range("K38").Value = 26000
$ws.Range("L38").Value = 39949.5
$ws.Range("M38").Value = -25584
$ws.Range("N38").Value = -40781.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 4374.5
$ws.Range("J15").Value = 4374.5
$ws.Range("L15").Value = 4374.5
$ws.Range("N15").Value = -4714.5
$ws.Range("H86").Value = 34665.668
$ws.Range("I86").Value = 34000
$ws.Range("K86").Value = 34000
$ws.Range("M86").Value = -32877
$ws.Range("H89").Value = 34665.668
$ws.Range("I89").Value = 34000
$ws.Range("K89").Value = 170000
$ws.Range("M89").Value = -164384
$ws.Range("H122").Value = 2023.5
$ws.Range("I122").Value = 1929.5
$ws.Range("J122").Value = 2399.5
$ws.Range("K122").Value = 5788.5
$ws.Range("L122").Value = 7198.5
$ws.Range("M122").Value = -3338.5
$ws.Range("N122").Value = -12098.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 4373.75
$ws.Range("J9").Value = 4427.143
$ws.Range("L9").Value = 13281.429
$ws.Range("N9").Value = -13729.429

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1321.2
$ws.Range("I102").Value = 1290.375
$ws.Range("J102").Value = 1650
$ws.Range("K102").Value = 1290.375
$ws.Range("L102").Value = 1650
$ws.Range("M102").Value = 331.625
$ws.Range("N102").Value = -4894

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3784.5715
$ws.Range("I7").Value = 3698.6
$ws.Range("K7").Value = 3698.6
$ws.Range("M7").Value = -3586.6
$ws.Range("H16").Value = 598
$ws.Range("I16").Value = 598
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 598
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -428
$ws.Range("N16").ClearContents()
$ws.Range("H40").Value = 52636520
$ws.Range("I40").Value = 83337090
$ws.Range("J40").Value = 6971.2856
$ws.Range("K40").Value = 83337090
$ws.Range("L40").Value = 6971.2856
$ws.Range("M40").Value = -83336954
$ws.Range("N40").Value = -7243.2856
$ws.Range("H61").Value = 3172.5881
$ws.Range("I61").Value = 3493.2856
$ws.Range("K61").Value = 3493.2856
$ws.Range("M61").Value = -3291.2856
$ws.Range("H93").Value = 1373.75
$ws.Range("I93").Value = 1098.4
$ws.Range("K93").Value = 1098.4
$ws.Range("M93").Value = 149.5999999999999
$ws.Range("H95").Value = 60000
$ws.Range("J95").Value = 60000
$ws.Range("L95").Value = 60000
$ws.Range("N95").Value = -65492
$ws.Range("H113").Value = 3172.5881
$ws.Range("I113").Value = 3493.2856
$ws.Range("K113").Value = 3493.2856
$ws.Range("M113").Value = -1323.2856
$ws.Range("H122").Value = 2229.3635
$ws.Range("I122").Value = 1902.7778
$ws.Range("K122").Value = 5708.3334
$ws.Range("M122").Value = -3258.3334
$ws.Range("H126").Value = 3784.5715
$ws.Range("I126").Value = 3698.6
$ws.Range("K126").Value = 11095.8
$ws.Range("M126").Value = -8625.799999999999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 9329
$ws.Range("J18").Value = 9493.5
$ws.Range("L18").Value = 9493.5
$ws.Range("N18").Value = -9839.5
$ws.Range("H62").Value = 7000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 7000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H97").Value = 21249
$ws.Range("J97").Value = 21249
$ws.Range("L97").Value = 21249
$ws.Range("N97").Value = -23231
$ws.Range("H136").Value = 5034.5
$ws.Range("I136").Value = 5983.76
$ws.Range("K136").Value = 17951.28
$ws.Range("M136").Value = -15401.28
